# Swap the contents of columns B:AC (columns 2-29) between paired rows,
# leaving column A (the sequential id) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(103,104),
    @(105,106),
    @(250,251),
    @(283,285),
    @(337,338),
    @(355,356),
    @(362,363),
    @(371,372),
    @(395,396),
    @(444,445),
    @(474,475),
    @(483,484),
    @(506,507),
    @(540,541),
    @(547,548),
    @(553,554),
    @(571,572)
)

$firstCol = 2   # column B
$lastCol  = 29  # column AC

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $row1Values = @{}
    $row2Values = @{}

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $row1Values[$c] = $ws.Cells.Item($r1, $c).Value2
        $row2Values[$c] = $ws.Cells.Item($r2, $c).Value2
    }

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r1, $c).Value2 = $row2Values[$c]
        $ws.Cells.Item($r2, $c).Value2 = $row1Values[$c]
    }
}
